$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("2021-10-16","overview","K02000001","United Kingdom",8404469,43423,148,138527),
  @("2021-10-17","overview","K02000001","United Kingdom",8449165,45140,57,138584),
  @("2021-10-18","overview","K02000001","United Kingdom",8497868,49156,45,138629),
  @("2021-10-19","overview","K02000001","United Kingdom",8541192,43738,223,138852),
  @("2021-10-20","overview","K02000001","United Kingdom",8589737,49139,179,139031),
  @("2021-10-21","overview","K02000001","United Kingdom",8641221,52009,115,139146),
  @("2021-10-22","overview","K02000001","United Kingdom",8689949,49298,180,139326),
  @("2021-10-23","overview","K02000001","United Kingdom",8734934,44985,135,139461),
  @("2021-10-24","overview","K02000001","United Kingdom",8773674,39962,72,139533)
)

$row = 431
foreach ($r in $data) {
  # Column A holds a date-like string ("YYYY-MM-DD"). Excel would normally
  # auto-convert that into a real date serial on assignment, so format the
  # cell as Text first, then strip the formatting back off afterwards so
  # the stored cell keeps the literal text with no leftover style index.
  $ws.Cells.Item($row, 1).NumberFormat = "@"
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 1).ClearFormats()

  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
  $ws.Cells.Item($row, 6).Value = $r[5]
  $ws.Cells.Item($row, 7).Value = $r[6]
  $ws.Cells.Item($row, 8).Value = $r[7]

  $row = $row + 1
}
